$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.469919769492789
$ws.Range("C2").Value = 5.402012787507591
$ws.Range("E2").Value = 16.45393243414715
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 29.56062815283716
$ws.Range("H2").Value = 14.37129007864003
$ws.Range("I2").Value = 20.43298986046508
$ws.Range("K2").Value = 9.423765959636974
$ws.Range("N2").Value = 17.98725513313039
$ws.Range("B3").Value = 9.147215534406818
$ws.Range("C3").Value = 5.09092083574448
$ws.Range("E3").Value = 15.52356817784021
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 29.53893829312201
$ws.Range("H3").Value = 14.41946763758516
$ws.Range("I3").Value = 20.51767391087926
$ws.Range("K3").Value = 9.198482505174914
$ws.Range("N3").Value = 18.05552134879807
$ws.Range("B4").Value = 8.945242350474301
$ws.Range("C4").Value = 4.88901637195436
$ws.Range("E4").Value = 14.92781344155641
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 29.53741982651899
$ws.Range("H4").Value = 14.45192169884533
$ws.Range("I4").Value = 20.5743314773648
$ws.Range("K4").Value = 9.059511626755471
$ws.Range("N4").Value = 18.09928469272128
$ws.Range("B5").Value = 8.862110661537438
$ws.Range("C5").Value = 4.804020669530364
$ws.Range("E5").Value = 14.67914688706191
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 29.53975760414241
$ws.Range("H5").Value = 14.46586755326404
$ws.Range("I5").Value = 20.5985880942536
$ws.Range("K5").Value = 9.00280526594771
$ws.Range("N5").Value = 18.1175848414952
$ws.Range("B6").Value = 8.848261062005806
$ws.Range("C6").Value = 4.789743772462155
$ws.Range("E6").Value = 14.63750911505343
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 29.54032398947682
$ws.Range("H6").Value = 14.46822671776668
$ws.Range("I6").Value = 20.60268632934189
$ws.Range("K6").Value = 8.993387338685006
$ws.Range("N6").Value = 18.12065177440406
$ws.Range("B7").Value = 8.94412436176391
$ws.Range("C7").Value = 4.887881063018356
$ws.Range("E7").Value = 14.92448329568454
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 29.53743939999179
$ws.Range("H7").Value = 14.45210686248218
$ws.Range("I7").Value = 20.57465388581218
$ws.Range("K7").Value = 9.058747048713693
$ws.Range("N7").Value = 18.09952960504512
$ws.Range("B8").Value = 9.359525504292037
$ws.Range("C8").Value = 5.297019938824136
$ws.Range("E8").Value = 16.13837621386226
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 29.55069529742658
$ws.Range("H8").Value = 14.38730434139746
$ws.Range("I8").Value = 20.46121913920082
$ws.Range("K8").Value = 9.346273973522115
$ws.Range("N8").Value = 18.01041084160255
$ws.Range("B9").Value = 10.13827044381388
$ws.Range("C9").Value = 6.01226047360421
$ws.Range("E9").Value = 18.38689457359149
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 29.67061920872602
$ws.Range("H9").Value = 14.28310342552385
$ws.Range("I9").Value = 20.27593340007377
$ws.Range("K9").Value = 9.90130602073336
$ws.Range("N9").Value = 17.85023320024073
$ws.Range("B10").Value = 10.68206710125324
$ws.Range("C10").Value = 6.48410012210165
$ws.Range("E10").Value = 20.02198744120049
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 29.81617748703473
$ws.Range("H10").Value = 14.22060474017653
$ws.Range("I10").Value = 20.16269703702666
$ws.Range("K10").Value = 10.29907818600158
$ws.Range("N10").Value = 17.74133305884141
$ws.Range("B11").Value = 10.92215236026342
$ws.Range("C11").Value = 6.687051736151436
$ws.Range("E11").Value = 20.72385057915873
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 29.89483567533492
$ws.Range("H11").Value = 14.19524854921563
$ws.Range("I11").Value = 20.11620547302218
$ws.Range("K11").Value = 10.4769725980814
$ws.Range("N11").Value = 17.69367529814501
$ws.Range("B12").Value = 11.01193665249546
$ws.Range("C12").Value = 6.762221562852197
$ws.Range("E12").Value = 20.98363225195871
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 29.92640198549851
$ws.Range("H12").Value = 14.18609089060453
$ws.Range("I12").Value = 20.09932653020493
$ws.Range("K12").Value = 10.54383015188674
$ws.Range("N12").Value = 17.67589735569162
$ws.Range("B13").Value = 10.99265149980556
$ws.Range("C13").Value = 6.746107290606929
$ws.Range("E13").Value = 20.92794988143827
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 29.91952462624062
$ws.Range("H13").Value = 14.18804336914099
$ws.Range("I13").Value = 20.10292933395909
$ws.Range("K13").Value = 10.52945476707707
$ws.Range("N13").Value = 17.67971421342498
$ws.Range("B14").Value = 10.92956206813413
$ws.Range("C14").Value = 6.693269751015183
$ws.Range("E14").Value = 20.74534300468651
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 29.89739704267327
$ws.Range("H14").Value = 14.19448622956051
$ws.Range("I14").Value = 20.11480224828335
$ws.Range("K14").Value = 10.48248354544755
$ws.Range("N14").Value = 17.6922073133809
$ws.Range("B15").Value = 10.89076842850628
$ws.Range("C15").Value = 6.66068594531014
$ws.Range("E15").Value = 20.63271080784093
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 29.88407475595089
$ws.Range("H15").Value = 14.19849057281323
$ws.Range("I15").Value = 20.12216948373912
$ws.Range("K15").Value = 10.45364433353481
$ws.Range("N15").Value = 17.69989468639041
$ws.Range("B16").Value = 10.66622299093918
$ws.Range("C16").Value = 6.47060125741612
$ws.Range("E16").Value = 19.97527778794497
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 29.81128667675246
$ws.Range("H16").Value = 14.22232387356013
$ws.Range("I16").Value = 20.1658367747876
$ws.Range("K16").Value = 10.28738497794513
$ws.Range("N16").Value = 17.74448533403654
$ws.Range("B17").Value = 10.52654256858247
$ws.Range("C17").Value = 6.350993830304373
$ws.Range("E17").Value = 19.56124375124772
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 29.76981515475792
$ws.Range("H17").Value = 14.23773381001916
$ws.Range("I17").Value = 20.19391424863791
$ws.Range("K17").Value = 10.18455797017895
$ws.Range("N17").Value = 17.77232104513505
$ws.Range("B18").Value = 10.44552013769755
$ws.Range("C18").Value = 6.281098432494926
$ws.Range("E18").Value = 19.31915191422912
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 29.7471338037697
$ws.Range("H18").Value = 14.24688649032124
$ws.Range("I18").Value = 20.21053592729778
$ws.Range("K18").Value = 10.12513086391709
$ws.Range("N18").Value = 17.78850859298717
$ws.Range("B19").Value = 10.41797288259489
$ws.Range("C19").Value = 6.257244169564154
$ws.Range("E19").Value = 19.23650426596693
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 29.73965579421802
$ws.Range("H19").Value = 14.2500350545303
$ws.Range("I19").Value = 20.21624471572246
$ws.Range("K19").Value = 10.10496329321886
$ws.Range("N19").Value = 17.79401989682134
$ws.Range("B20").Value = 10.54148299491305
$ws.Range("C20").Value = 6.363840153116667
$ws.Range("E20").Value = 19.60572699042331
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 29.77410862892435
$ws.Range("H20").Value = 14.23606344043787
$ws.Range("I20").Value = 20.19087644395488
$ws.Range("K20").Value = 10.19553395494043
$ws.Range("N20").Value = 17.76933955948771
$ws.Range("B21").Value = 10.94812426129617
$ws.Range("C21").Value = 6.708835118154574
$ws.Range("E21").Value = 20.79914160683686
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 29.90384822786607
$ws.Range("H21").Value = 14.19258173429729
$ws.Range("I21").Value = 20.1112951353281
$ws.Range("K21").Value = 10.49629441341239
$ws.Range("N21").Value = 17.6885304982289
$ws.Range("B22").Value = 11.20726523608607
$ws.Range("C22").Value = 6.924496425083227
$ws.Range("E22").Value = 21.54415763520837
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 29.99900993824632
$ws.Range("H22").Value = 14.16675387122167
$ws.Range("I22").Value = 20.06352026633907
$ws.Range("K22").Value = 10.68987452232216
$ws.Range("N22").Value = 17.63728443632416
$ws.Range("B23").Value = 11.06958753703318
$ws.Range("C23").Value = 6.810292133217727
$ws.Range("E23").Value = 21.14971541479037
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 29.94727548572479
$ws.Range("H23").Value = 14.18030104430328
$ws.Range("I23").Value = 20.0886295235368
$ws.Range("K23").Value = 10.58685136361691
$ws.Range("N23").Value = 17.66449251836815
$ws.Range("B24").Value = 10.53473065910903
$ws.Range("C24").Value = 6.358035854834708
$ws.Range("E24").Value = 19.58562876150476
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 29.77216393141411
$ws.Range("H24").Value = 14.23681770097839
$ws.Range("I24").Value = 20.19224834186392
$ws.Range("K24").Value = 10.19057267406676
$ws.Range("N24").Value = 17.77068691448661
$ws.Range("B25").Value = 9.932152159568167
$ws.Range("C25").Value = 5.828163280400064
$ws.Range("E25").Value = 17.74906723250236
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 29.62808726584334
$ws.Range("H25").Value = 14.30883183695821
$ws.Range("I25").Value = 20.32205519575081
$ws.Range("K25").Value = 9.752594785740209
$ws.Range("N25").Value = 17.89201541422073
